# Insert a new data row at row 14 (pushes existing rows 14-62 down to 15-63)
# and populate it with a new price observation for "Vega Modelo de Temuco" /
# Arándano (blue), matching the other rows' shared columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(14).Insert()

$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44525
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100101
$ws.Range("H14").Value = "Berries"
$ws.Range("I14").Value = 100101001
$ws.Range("J14").Value = "Arándano (blue)"
$ws.Range("K14").Value = "Sin especificar"
$ws.Range("L14").Value = "Primera"
$ws.Range("M14").Value = 1450
$ws.Range("N14").Value = 3000
$ws.Range("O14").Value = 3200
$ws.Range("P14").Value = 3110
$ws.Range("Q14").Value = "$/kilo"
$ws.Range("R14").Value = "Región del Maule"
$ws.Range("S14").Value = 3110
$ws.Range("T14").Value = 1
